# Remove the <w:contextualSpacing w:val="0"/> element from every paragraph's
# properties (w:pPr). The Word object model exposed by this runtime does not
# surface a ContextualSpacing property on Paragraph/ParagraphFormat, so we
# fall back to round-tripping each paragraph through WordOpenXML/InsertXML:
# read the paragraph's own OOXML, strip the one element, and write it back
# into the same range.

$d = $word.ActiveDocument

# InsertXML re-numbers relationship ids (r:id) using the scoped package it
# was handed, which does not line up with the real document.xml.rels ids -
# that silently corrupts hyperlinks touched by the round trip. Save every
# hyperlink address up front and restore it afterwards.
$savedHyperlinkAddresses = @()
for ($h = 1; $h -le $d.Hyperlinks.Count; $h++) {
    $savedHyperlinkAddresses += $d.Hyperlinks($h).Address
}

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $rng = $para.Range

    $wx = $rng.WordOpenXML
    if ($wx -eq $null) { continue }
    if (-not $wx.Contains("contextualSpacing")) { continue }

    $bodyTag = "<w:body>"
    $bodyStart = $wx.IndexOf($bodyTag)
    if ($bodyStart -lt 0) { continue }
    $bodyStart = $bodyStart + $bodyTag.Length
    $bodyEnd = $wx.IndexOf("</w:body>", $bodyStart)
    if ($bodyEnd -lt 0) { continue }
    $bodyInner = $wx.Substring($bodyStart, $bodyEnd - $bodyStart)

    # Isolate just this paragraph's own <w:p ...>...</w:p> node - the range's
    # WordOpenXML can include a trailing synthesized empty paragraph / sectPr
    # after it, which must not be reinserted.
    $pStart = $bodyInner.IndexOf("<w:p ")
    if ($pStart -lt 0) { $pStart = $bodyInner.IndexOf("<w:p>") }
    if ($pStart -lt 0) { continue }
    $closeTag = "</w:p>"
    $closeIdx = $bodyInner.IndexOf($closeTag, $pStart)
    if ($closeIdx -lt 0) { continue }
    $pEnd = $closeIdx + $closeTag.Length
    $firstPara = $bodyInner.Substring($pStart, $pEnd - $pStart)

    if (-not $firstPara.Contains("contextualSpacing")) { continue }

    $cleaned = $firstPara -replace '<w:contextualSpacing[^>]*/>', ''

    $rng.InsertXML($cleaned)
}

# Restore any hyperlink addresses InsertXML mangled above.
for ($h = 1; $h -le $d.Hyperlinks.Count; $h++) {
    $wanted = $savedHyperlinkAddresses[$h - 1]
    if ($d.Hyperlinks($h).Address -ne $wanted) {
        $d.Hyperlinks($h).Address = $wanted
    }
}
